# Commit: "Test copy data from manmos to excel"
#
# - Adds a new "collectionDataError" mail-template row (A6/B6) to the
#   "Mail" sheet, mirroring the existing error-mail rows above it.
# - Makes "Mail" the active/selected sheet (it previously was "Path"),
#   with B6 selected as the active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mail")

# New email-template body for the "collection data" timeout error.
$body = @'
{
"sendTo":"{sendTo}",
"subject":"{robotName} - Failed to load collection data window.",
"body":
"<html>
<head>
<style>
div {font-family:Cordia New, Arial, sans-serif; font-size:24px;}
</style>
</head>
<body>
<div>
Dear all,<br/><br/>
Robot couldn't load collection data from manmos.<br/>
Because manmos working slowly more than 5 minute<br/><br/>
Thank you<br/><br/>
</div>
<div>
<b>Robotic Process Automation (RPA)<br/>
APP Department - Head Office<br/>
<img src='https://www.jobtni.com/files/company/logo/281.jpg' alt='' height='80'><br/>
NHK Spring (Thailand) CO.,LTD.<br/>
Phone : (+66)2-730-2200 Ext: 2433<br/>
E-mail : </b><a href='mailto:chakrit.pok@nhkspg.co.th?Subject=Contact Problem' target='_top'>chakrit.pok@nhkspg.co.th</a>
</div>
</body>
</html>"
}
'@

# Activate "Mail" first (it becomes the workbook's active/selected tab,
# replacing "Path") before writing the new row's data.
$ws.Activate()

$ws.Range("A6").Value = "collectionDataError"
$ws.Range("B6").Value = $body

# This row's wrapped text needs the maximum row height, same as the two
# other big HTML-template rows above it (A4:B4 / A5:B5).
$ws.Rows.Item(6).RowHeight = 409.5

# Final selection lands on B6, scrolled so row 6 is visible.
$ws.Range("B6").Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
